# Regenerate s_vals data to filter save games.
# Updates the numeric stat columns (B:E and G) for rows 2-5.
# Column F (Win flag) is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.04763786555579896, 0.3127903958511391, 0.8054896365839992, 8.660232485948974, 9.826150383939911)
    3 = @(0.3048080303191223, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.274871460341982)
    4 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    5 = @(3.230985683306322, 1.667794583268128, 9844.520545567508, 645.3272768299601, 10494.74660266404)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - sum
}
